$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("wgs")
$ws.Activate()

# Insert 3 new columns before column B (old B..AA shift to E..AD)
$ws.Range("B1:D1").EntireColumn.Insert()

# Populate the new header cells
$ws.Range("B1").Value2 = "Experiment Alias"
$ws.Range("C1").Value2 = "Project"
$ws.Range("D1").Value2 = "Secondary Project"
